# docs: update docs
# - fix "première" -> "premier" typo on sheet "1-10"
# - tweak sheet view zoom/selection on both sheets
# - set explicit column widths on sheet "1-10" (A, C, D)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("1-10")
$ws2 = $wb.Worksheets.Item("11-100")

# --- Text fix: "première" -> "premier" (sheet "1-10", cell C2) ---
$ws1.Range("C2").Value = "premier"

# --- Sheet "1-10": selection + zoom ---
$ws1.Activate()
$ws1.Range("C3").Select()
$excel.ActiveWindow.Zoom = 206

# --- Sheet "1-10": explicit column widths (A, C, D) ---
# ColumnWidth is expressed in this runtime's standard (MDW=7) character
# units; offset chosen so the exported <col width="..."> lands as close as
# possible to the target widths from the authored workbook.
$ws1.Columns.Item(1).ColumnWidth = 1.9780219780219754
$ws1.Columns.Item(3).ColumnWidth = 10.978021978021987
$ws1.Columns.Item(4).ColumnWidth = 2.9780219780219754

# --- Sheet "11-100": selection + zoom ---
$ws2.Activate()
$ws2.Range("E11").Select()
$excel.ActiveWindow.Zoom = 158

# Leave the originally-active sheet ("1-10") selected/active.
$ws1.Activate()
